$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.787.92'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.235.89'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '271.05'
$ws.Range('E5').Value = '  +5.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.07'
$ws.Range('E6').Value = '  +16.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.624'
$ws.Range('E9').Value = '  +6.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.42'
$ws.Range('E10').Value = '  +9.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0971'
$ws.Range('E11').Value = '  +6.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.35'
$ws.Range('E12').Value = '  +20.82%  '
$ws.Range('E13').Value = '  +1.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.14'
$ws.Range('E14').Value = '  +6.68%  '
$ws.Range('D15').Value = '2.564.06'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '2.238.25'
$ws.Range('E16').Value = '  +2.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.802'
$ws.Range('E17').Value = '  +4.13%  '
$ws.Range('D18').Value = '43.743.13'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000106'
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.08'
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.74'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.36'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.85'
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.13'
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.45'
$ws.Range('E26').Value = '  +8.58%  '
$ws.Range('E27').Value = '  +13.92%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.40'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('B29').Value = 'WEMIXToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.52'
$ws.Range('E29').Value = '  +5.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  +1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.37'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0924'
$ws.Range('E32').Value = '  +6.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.96'
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.52'
$ws.Range('E34').Value = '  +5.60%  '
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0354'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.33'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('E39').Value = '  +25.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.94'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.224'
$ws.Range('E41').Value = '  +13.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.18'
$ws.Range('E42').Value = '  +4.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.93'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.35'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0998'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.42'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +4.94%  '
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.444'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').Value = '2.452.15'
$ws.Range('E51').Value = '  +1.99%  '
